# Auto-generated edit script updating cryptos price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "42.388.28"
$ws.Cells.Item(2, 5).Value = "  +0.61%  "
$ws.Cells.Item(3, 4).Value = "2.246.52"
$ws.Cells.Item(3, 5).Value = "  +0.28%  "
$ws.Cells.Item(4, 5).Value = "  +0.07%  "
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = "246.14"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.30%  "
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.631"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +0.17%  "
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = "75.71"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.82%  "
$ws.Cells.Item(8, 5).Value = "  +0.07%  "
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.619"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  -1.72%  "
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = "44.13"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +8.45%  "
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0948"
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.02%  "
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = "7.19"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +0.46%  "
$ws.Cells.Item(13, 5).Value = "  -1.65%  "
$ws.Cells.Item(14, 2).Value = "Chainlink"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = "14.58"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.93%  "
$ws.Cells.Item(15, 2).Value = "Polygon"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.858"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.41%  "
$ws.Cells.Item(16, 2).Value = "WrappedEther"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Cells.Item(16, 4).Value = "2.233.91"
$ws.Cells.Item(16, 5).Value = "  -0.42%  "
$ws.Cells.Item(17, 2).Value = "WrappedBTC"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(17, 4).Value = "42.224.16"
$ws.Cells.Item(17, 5).Value = "  +0.55%  "
$ws.Cells.Item(18, 2).Value = "ShibaInu"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0000102"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +4.12%  "
$ws.Cells.Item(19, 2).Value = "Uniswap"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = "6.17"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.18%  "
$ws.Cells.Item(20, 2).Value = "Litecoin"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = "72.19"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.03%  "
$ws.Cells.Item(21, 2).Value = "ImmutableX"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.25"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +2.84%  "
$ws.Cells.Item(22, 2).Value = "BitcoinCash"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = "231.80"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.20%  "
$ws.Cells.Item(23, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = "9.05"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +26.21%  "
$ws.Cells.Item(24, 2).Value = "Dai"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.07%  "
$ws.Cells.Item(25, 2).Value = "Cosmos"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = "11.47"
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +2.54%  "
$ws.Cells.Item(26, 2).Value = "WEMIXToken"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = "3.62"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -2.85%  "
$ws.Cells.Item(27, 2).Value = "PancakeSwap"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.31"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -0.03%  "
$ws.Cells.Item(28, 2).Value = "Toncoin"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.20"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.79%  "
$ws.Cells.Item(29, 2).Value = "Monero"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = "167.60"
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.64%  "
$ws.Cells.Item(30, 2).Value = "EthereumClassic"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = "20.67"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.81%  "
$ws.Cells.Item(31, 2).Value = "Hedera"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0825"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -3.22%  "
$ws.Cells.Item(32, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = "30.99"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -4.91%  "
$ws.Cells.Item(33, 5).Value = "  +1.13%  "
$ws.Cells.Item(34, 2).Value = "Filecoin"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.34"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +9.37%  "
$ws.Cells.Item(35, 2).Value = "Stellar"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.126"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -0.01%  "
$ws.Cells.Item(36, 2).Value = "RenderToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = "4.52"
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -0.30%  "
$ws.Cells.Item(37, 2).Value = "VeChain"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.0315"
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +6.00%  "
$ws.Cells.Item(38, 2).Value = "Celestia"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = "13.95"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +5.27%  "
$ws.Cells.Item(39, 2).Value = "LidoDAOToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.18"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.76%  "
$ws.Cells.Item(40, 2).Value = "THORChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = "5.80"
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -1.59%  "
$ws.Cells.Item(41, 2).Value = "MultiversX"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = "64.03"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +6.57%  "
$ws.Cells.Item(42, 2).Value = "Algorand"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.202"
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.46%  "
$ws.Cells.Item(43, 2).Value = "Aave"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = "107.67"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -7.39%  "
$ws.Cells.Item(44, 2).Value = "FraxShare"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = "8.81"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +0.70%  "
$ws.Cells.Item(45, 2).Value = "Cronos"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +2.53%  "
$ws.Cells.Item(46, 2).Value = "BinanceUSD"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.996"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.09%  "
$ws.Cells.Item(47, 2).Value = "ARBITRUM"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.13"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.10%  "
$ws.Cells.Item(48, 2).Value = "TrustWalletToken"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = "1.19"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +1.98%  "
$ws.Cells.Item(49, 2).Value = "NEARProtocol"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.35"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.90%  "
$ws.Cells.Item(50, 2).Value = "HuobiToken"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = "2.71"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.01%  "
$ws.Cells.Item(51, 2).Value = "WOONetwork"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = "0.424"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +16.29%  "
